$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows appended to the master-template_type table, mirroring the
# existing "reg-ack-template-partN" rows (eng / ara / fra variants).
$ws.Range("A122").Value = "reg-ack-template-part4"
$ws.Range("B122").Value = "Registration Acknowledgement Template - Part 4"
$ws.Range("C122").Value = "eng"
$ws.Range("D122").Value = $true
$ws.Range("E122").Value = "superadmin"
$ws.Range("F122").Value = "now()"

$ws.Range("A123").Value = "reg-ack-template-part4"
$ws.Range("B123").Value = "نموذج شكر التسجيل"
$ws.Range("C123").Value = "ara"
$ws.Range("D123").Value = $true
$ws.Range("E123").Value = "superadmin"
$ws.Range("F123").Value = "now()"

$ws.Range("A124").Value = "reg-ack-template-part4"
$ws.Range("B124").Value = "accusé de réception"
$ws.Range("C124").Value = "fra"
$ws.Range("D124").Value = $true
$ws.Range("E124").Value = "superadmin"
$ws.Range("F124").Value = "now()"

# Match the post-edit selection left behind in the saved file (user had
# selected from the next empty row to the bottom of the sheet).
$r1 = $ws.Rows.Item(125)
$r2 = $ws.Rows.Item(1048576)
$ws.Range($r1, $r2).Select() | Out-Null
